$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (Treatment) from 18 to 22 character-units
$ws.Columns.Item(1).ColumnWidth = 21.1

# Row 62
$ws.Range("A62").Value = "Transmission correct"
$ws.Range("B62").Value = "yp80d4yg"
$ws.Range("C62").Value = "Training phase"
$ws.Range("D62").Value = 1
$ws.Range("E62").Value = "['Purple', 'Orange', 'Green']"
$ws.Range("F62").Value = "[['Red', ''], ['Red', ''], ['Red', '']]"
$ws.Range("G62").Value = "[None, None, None]"
$ws.Range("H62").Value = "['2', '2', '2']"

# Row 63
$ws.Range("A63").Value = "Transmission correct"
$ws.Range("B63").Value = "yp80d4yg"
$ws.Range("C63").Value = "Training phase"
$ws.Range("D63").Value = 2
$ws.Range("E63").Value = "['Green', 'Green', 'Orange']"
$ws.Range("F63").Value = "[['Red', ''], ['Red', ''], ['Red', '']]"
$ws.Range("G63").Value = "[None, None, None]"
$ws.Range("H63").Value = "['2', '2', '2']"

# Row 64
$ws.Range("A64").Value = "Transmission correct"
$ws.Range("B64").Value = "yp80d4yg"
$ws.Range("C64").Value = "Training phase"
$ws.Range("D64").Value = 3
$ws.Range("E64").Value = "['Purple', 'Green', 'Purple']"
$ws.Range("F64").Value = "[['Red', ''], ['Red', ''], ['Red', '']]"
$ws.Range("G64").Value = "[None, None, None]"
$ws.Range("H64").Value = "['2', '2', '2']"

# Row 65
$ws.Range("A65").Value = "Transmission correct"
$ws.Range("B65").Value = "yp80d4yg"
$ws.Range("C65").Value = "Training phase"
$ws.Range("D65").Value = 4
$ws.Range("E65").Value = "['Purple', 'Orange', 'Orange']"
$ws.Range("F65").Value = "[['Red', ''], ['Red', ''], ['Red', '']]"
$ws.Range("G65").Value = "[None, None, None]"
$ws.Range("H65").Value = "['2', '2', '2']"

# Row 66
$ws.Range("A66").Value = "Transmission correct"
$ws.Range("B66").Value = "yp80d4yg"
$ws.Range("C66").Value = "Training phase"
$ws.Range("D66").Value = 5
$ws.Range("E66").Value = "['Green', 'Orange', 'Purple']"
$ws.Range("F66").Value = "[['Red', ''], ['Red', ''], ['Red', '']]"
$ws.Range("G66").Value = "[None, None, None]"
$ws.Range("H66").Value = "['2', '2', '2']"

# Row 67
$ws.Range("A67").Value = "Transmission correct"
$ws.Range("B67").Value = "yp80d4yg"
$ws.Range("C67").Value = "Test 1"
$ws.Range("D67").Value = 1
$ws.Range("E67").Value = "['Green', 'Yellow', 'Purple', 'Red', 'Orange', 'Blue']"
$ws.Range("F67").Value = "[['Blue', 'Blue'], ['Blue', 'Blue'], ['Blue', 'Blue'], ['Blue', 'Blue'], ['Blue', 'Blue'], ['Blue', 'Blue']]"
$ws.Range("G67").Value = "[None, None, None, None, None, None]"
$ws.Range("H67").Value = "['0', '0', '0', '0', '0', '0']"

# Row 68
$ws.Range("A68").Value = "Transmission correct"
$ws.Range("B68").Value = "yp80d4yg"
$ws.Range("C68").Value = "Test 1"
$ws.Range("D68").Value = 1
$ws.Range("E68").Value = "['Green', 'Yellow', 'Purple', 'Red', 'Orange', 'Blue']"
$ws.Range("F68").Value = "[['Blue', 'Yellow'], ['Yellow', 'Yellow'], ['Red', 'Blue'], ['Red', 'Red'], ['Yellow', 'Red'], ['Blue', 'Blue']]"
$ws.Range("G68").Value = "[None, None, None, None, None, None]"
$ws.Range("H68").Value = "['0', '0', '0', '0', '0', '0']"
$ws.Range("I68").NumberFormat = "@"
$ws.Range("I68").Value = "0.30"
$ws.Range("I68").Style = "Normal"
